$d = $word.ActiveDocument

function Replace-Exact($OldText, $NewText, $SetBold) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Text = $OldText
    $find.Replacement.ClearFormatting()
    if ($SetBold) {
        $find.Replacement.Font.Bold = 1
    }
    $find.Replacement.Text = $NewText
    $result = $find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $true, $NewText, 2)
    return $result
}

# 1. Title line - bold on + text change
Replace-Exact "Documento de responsabilidades principais da equipe de design" `
              "Documento de responsabilidade principal da equipe de projeto" `
              $true

# 2. Intro sentence - capitalization change
Replace-Exact ": Este documento descreve as principais responsabilidades de todos os membros da equipe de design do Instituto de Design Gráfico." `
              ": este documento descreve as principais responsabilidades de todos os membros da equipe de design do Instituto de Design Gráfico." `
              $false

# 3. Collaboration bullet
Replace-Exact ": Trabalhe de forma colaborativa com outros designers, desenvolvedores e partes interessadas para criar designs de alta qualidade que atendam aos requisitos do projeto." `
              ": trabalhe em colaboração com outros designers, desenvolvedores e stakeholders para criar designs de alta qualidade que atendam aos requisitos do projeto." `
              $false

# 4. Visual design bullet
Replace-Exact ": crie designs visualmente atraentes que sejam fáceis de usar, acessíveis e responsivos." `
              ":  crie designs visualmente atraentes que sejam acessíveis, responsivos e fáceis de usar pelo usuário." `
              $false

# 5. Communication bullet
Replace-Exact ": Comunique-se efetivamente com os membros da equipe, partes interessadas e clientes para garantir que os requisitos do projeto sejam atendidos." `
              ": comunique-se de forma efetiva com os membros da equipe, stakeholders e clientes para garantir que os requisitos do projeto sejam atendidos." `
              $false

# 6. Research bullet
Replace-Exact ": Realizar pesquisas para identificar necessidades, preferências e comportamentos do usuário para informar as decisões de design." `
              ": faça pesquisas para identificar as necessidades, preferências e comportamentos do usuário para informar as decisões de design." `
              $false

# 7. Usability testing bullet
Replace-Exact ": conduza testes de usabilidade para garantir que os designs atendam às necessidades do usuário e sejam acessíveis a todos os usuários." `
              ": faça testes de usabilidade para garantir que os designs atendam às necessidades dos usuários e sejam acessíveis para todos os usuários." `
              $false

# 8. Documentation bullet
Replace-Exact ": Crie e mantenha a documentação de projeto, incluindo especificações de projeto, guias de estilo e padrões de design." `
              ": crie e mantenha a documentação do design, incluindo especificações do design, guias de estilo e padrões de design." `
              $false

# 9. "Desenvolvimento" -> "Desenvolvimento profissional" (bold run, keep bold)
Replace-Exact "Desenvolvimento" `
              "Desenvolvimento profissional" `
              $false

# 10. Professional development detail sentence
Replace-Exact " Profissional: Mantenha-se atualizado com as últimas tendências, ferramentas e tecnologias de design para melhorar a qualidade e a eficiência do projeto." `
              ": esteja a par das últimas tendências, ferramentas e tecnologias de design para melhorar a qualidade e a eficiência do design." `
              $false

# 11. Leadership bullet
Replace-Exact " Liderar a equipe de design e fornecer orientação aos designers juniores." `
              " lidere a equipe de design e forneça orientação aos designers juniores." `
              $false
